$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.3418048620224
$ws.Range("B1").Value = 2.656242370605469
$ws.Range("D1").Value = 1.559033393859863
$ws.Range("E1").Value = 0.92879319190979
